# add c_click with many issues
# Update the "수량" (quantity) counts in the "기타" (misc) sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("기타")

$ws.Range("C3").Value = 1
$ws.Range("C11").Value = 5
$ws.Range("C12").Value = 7
$ws.Range("C14").Value = 20
